$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure affected cells retain their literal text representation
# (values in this sheet are stored as text, not numbers/percentages).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "245.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.46%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "28.37"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "-3.54%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.288"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "1.97%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05711"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-0.41%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.639"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.17%"

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.41%"

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8545"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.46%"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8909"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "2.63%"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1392"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "1.83%"

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07087"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.03%"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "4.95%"

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09229"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.66%"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001527"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.56%"

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0005999"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-0.23%"

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005881"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-1.44%"

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.497"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.08%"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.48%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.57%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.03339"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.90%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.1307"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "1.25%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.480"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-0.45%"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.98%"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.14%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001219"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-0.57%"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-16.86%"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-0.90%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03798"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.17%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1066"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "-0.55%"

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-35.25%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002199"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.69%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009436"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-0.46%"

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005282"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.04%"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.07%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.08909"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "56.24%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.002258"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.04%"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.07%"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.07%"
